$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Time Log-Quadrants")

$ws.Range("F9").Value = "Sleep"
$ws.Range("F10").Value = "Sleep"
$ws.Range("F11").Value = "Sleep"
$ws.Range("F12").Value = "Sleep"
$ws.Range("F13").Value = "Sleep"
$ws.Range("F14").Value = "Sleep"
$ws.Range("F15").Value = "Sleep"
$ws.Range("F16").Value = "Sleep"
$ws.Range("C17").Value = "Eat"
$ws.Range("D17").Value = "Eat"
$ws.Range("E17").Value = "Eat"
$ws.Range("F17").Value = "Shower"
$ws.Range("F18").Value = "Eat"
$ws.Range("B22").Value = "Eat"
$ws.Range("B23").Value = "Eat"
$ws.Range("E24").Value = "Homework"
$ws.Range("E25").Value = "Homework"
$ws.Range("F25").Value = "Homework"
$ws.Range("D26").Value = "Eat"
$ws.Range("E26").Value = "Homework"
$ws.Range("F26").Value = "Homework"
$ws.Range("D27").Value = "Eat"
$ws.Range("E27").Value = "Eat"
$ws.Range("F27").Value = "Homework"
$ws.Range("E28").Value = "Homework"
$ws.Range("F28").Value = "Homework"
$ws.Range("E29").Value = "Homework"
$ws.Range("F29").Value = "Homework"
$ws.Range("E30").Value = "Homework"
$ws.Range("F30").Value = "Homework"
$ws.Range("E31").Value = "Homework"
$ws.Range("F31").Value = "Homework"
$ws.Range("E32").Value = "Homework"
$ws.Range("F32").Value = "Homework"
$ws.Range("E33").Value = "Homework"
$ws.Range("F33").Value = "Homework"
$ws.Range("E34").Value = "Homework"
$ws.Range("F34").Value = "Homework"
$ws.Range("J34").Value = "Eat"
$ws.Range("E35").Value = "Homework"
$ws.Range("F35").Value = "Homework"
$ws.Range("E36").Value = "Homework"
$ws.Range("E37").Value = "Homework"
$ws.Range("E38").Value = "Homework"
$ws.Range("B39").Value = "Eat"
$ws.Range("E39").Value = "Homework"
$ws.Range("B40").Value = "Eat"
$ws.Range("E40").Value = "Homework"
$ws.Range("J40").Value = "Class"
$ws.Range("C45").Value = "Eat"
$ws.Range("F45").Value = "Work"
$ws.Range("C46").Value = "Eat"
$ws.Range("F46").Value = "Work"
$ws.Range("F47").Value = "Work"
$ws.Range("F48").Value = "Work"
$ws.Range("E55").Value = "Eat"
$ws.Range("F55").Value = "Clean apartment"
$ws.Range("E56").Value = "Clean apartment"
$ws.Range("F56").Value = "Shower"
$ws.Range("D57").Value = "Eat"
$ws.Range("E57").Value = "Clean apartment"
$ws.Range("F57").Value = "Eat"
$ws.Range("D58").Value = "Eat"
$ws.Range("E58").Value = "Clean apartment"
$ws.Range("F58").Value = "Homework"
$ws.Range("C59").Value = "Eat"
$ws.Range("E59").Value = "Clean apartment"
$ws.Range("F59").Value = "Homework"
$ws.Range("E60").Value = "Homework"
$ws.Range("F60").Value = "Homework"
$ws.Range("E61").Value = "Homework"
$ws.Range("F61").Value = "Homework"
$ws.Range("E62").Value = "Homework"
$ws.Range("F62").Value = "Homework"
$ws.Range("E63").Value = "Homework"
$ws.Range("F63").Value = "Homework"
$ws.Range("E64").Value = "Homework"
$ws.Range("F64").Value = "Homework"
$ws.Range("E65").Value = "Sleep"
$ws.Range("F65").Value = "Homework"
$ws.Range("B66").Value = "Eat"
$ws.Range("E66").Value = "Sleep"
$ws.Range("F66").Value = "Homework"
$ws.Range("E67").Value = "Sleep"
$ws.Range("F67").Value = "Homework"
$ws.Range("E68").Value = "Sleep"
$ws.Range("F68").Value = "Homework"
$ws.Range("E69").Value = "Sleep"
$ws.Range("F69").Value = "Homework"
$ws.Range("E70").Value = "Sleep"
$ws.Range("F70").Value = "Homework"
$ws.Range("E71").Value = "Sleep"
$ws.Range("F71").Value = "Homework"
$ws.Range("E72").Value = "Sleep"
$ws.Range("F72").Value = "Homework"
$ws.Range("E73").Value = "Sleep"
$ws.Range("F73").Value = "Homework"
$ws.Range("E74").Value = "Sleep"
$ws.Range("F74").Value = "Homework"
$ws.Range("E75").Value = "Sleep"
$ws.Range("F75").Value = "Homework"
$ws.Range("E76").Value = "Sleep"
$ws.Range("F76").Value = "Homework"
$ws.Range("E77").Value = "Sleep"
$ws.Range("F77").Value = "Homework"
$ws.Range("E78").Value = "Sleep"
$ws.Range("F78").Value = "Homework"
$ws.Range("E79").Value = "Sleep"
$ws.Range("F79").Value = "Homework"
$ws.Range("E80").Value = "Sleep"
$ws.Range("F80").Value = "Homework"

# Update view state: active cell G26, scroll to top-left A1
$ws.Activate()
$ws.Range("G26").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
